# feat(data): update data source
# Adds a 4th episode to each of the three seasons on the "episodios" sheet
# (with a blank separator row between seasons, matching the existing
# layout convention), and refreshes the active selections on both sheets.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("episodios")

# --- Season 1: insert "episodio-4_temporada-1" + blank separator after row 4 ---
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(6).Insert()

$ws.Cells.Item(5,1).Value = "episodio-4_temporada-1"
$ws.Cells.Item(5,2).Value = "Temproada_1"
$ws.Cells.Item(5,3).Value = "Título Episodio 4 (Temporada 1)"
$ws.Cells.Item(5,4).Value = 4
$ws.Cells.Item(5,5).Value = "Descripción Episodio 4 Temporada 1"
$ws.Cells.Item(5,6).Value = "https://test-videos.co.uk/vids/bigbuckbunny/mp4/h264/720/Big_Buck_Bunny_720_10s_2MB.mp4"
$ws.Cells.Item(5,7).Value = "Episodio-3_Temporada-1.png"
$ws.Cells.Item(5,8).Value = 2022
$ws.Cells.Item(5,9).Value = "50min"
$ws.Cells.Item(5,10).Value = "Temporada_1.png"

# --- Season 2: insert "episodio-4_temporada-2" + blank separator after the
#     last season-2 row (originally row 7, now row 9 after the season-1 insert) ---
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(11).Insert()

$ws.Cells.Item(10,1).Value = "episodio-4_temporada-2"
$ws.Cells.Item(10,2).Value = "Temporada_2"
$ws.Cells.Item(10,3).Value = "Título Episodio 4 (Temporada 2)"
$ws.Cells.Item(10,4).Value = 4
$ws.Cells.Item(10,5).Value = "Descripción Episodio 4 Temporada 2"
$ws.Cells.Item(10,6).Value = "https://test-videos.co.uk/vids/bigbuckbunny/mp4/h264/720/Big_Buck_Bunny_720_10s_2MB.mp4"
$ws.Cells.Item(10,7).Value = "Episodio-3_Temporada-2.png"
$ws.Cells.Item(10,8).Value = 2022
$ws.Cells.Item(10,9).Value = "50min"
$ws.Cells.Item(10,10).Value = "Temporada_1.png"

# --- Season 3: append "episodio-4_temporada-3" as a new row 15 at the end ---
$ws.Rows.Item(15).Insert()

$ws.Cells.Item(15,1).Value = "episodio-4_temporada-3"
$ws.Cells.Item(15,2).Value = "Temporada_3"
$ws.Cells.Item(15,3).Value = "Título Episodio 4 (Temporada 3)"
$ws.Cells.Item(15,4).Value = 4
$ws.Cells.Item(15,5).Value = "Descripción Episodio 4 Temporada 3"
$ws.Cells.Item(15,6).Value = "https://test-videos.co.uk/vids/bigbuckbunny/mp4/h264/720/Big_Buck_Bunny_720_10s_2MB.mp4"
$ws.Cells.Item(15,7).Value = "Episodio-3_Temporada-3.png"
$ws.Cells.Item(15,8).Value = 2022
$ws.Cells.Item(15,9).Value = "50min"
$ws.Cells.Item(15,10).Value = "Temporada_1.png"

# The blank separator rows should only carry formatting in A:F (matching the
# sheet's existing spacer-row convention); fully clear the G:J cells the
# row-insert dragged along so they don't linger as empty styled cells.
$ws.Range("G6:J6").Clear()
$ws.Range("G11:J11").Clear()

# The original sheet's trailing placeholder cell (old "G11", empty but
# styled) rode along with the inserts and is now sitting on row 16, past
# the new table. Clear it so the used range ends cleanly at row 15.
$ws.Range("A16:J16").Clear()

# --- Refresh view selections to match the post-edit cursor position ---
$ws.Activate()
$ws.Range("B16").Select()

$ws1 = $wb.Worksheets.Item("temporadas")
$ws1.Activate()
$ws1.Range("E2").Select()

# Leave "episodios" as the active/selected sheet, matching the source file.
$ws.Activate()

Write-Output "done"
